$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point is inserted at row 158, pushing all subsequent
# rows (old 158..269) down by one (to 159..270).
$ws.Rows("158:158").Insert()

$ws.Cells.Item(158, 1).Value = 4
$ws.Cells.Item(158, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(158, 3).Value = "Los Lagos"
$ws.Cells.Item(158, 4).Value = 44957
$ws.Cells.Item(158, 5).Value = 10
$ws.Cells.Item(158, 6).Value = 100112009
$ws.Cells.Item(158, 7).Value = "Acelga"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 80
$ws.Cells.Item(158, 11).Value = 10000
$ws.Cells.Item(158, 12).Value = 10000
$ws.Cells.Item(158, 13).Value = 10000
$ws.Cells.Item(158, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(158, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(158, 16).Value = 833
$ws.Cells.Item(158, 17).Value = 12
$ws.Cells.Item(158, 18).Value = "Hortaliza"
